$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cell F1 ("time_taken") matching the style of the other header cells
$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

# Populate time_taken column for data rows 2-81
$ws.Cells.Item(2, 6).Value = "2021-10-05 13:39:34.561569"
$ws.Cells.Item(3, 6).Value = "2021-10-05 13:39:34.561582"
$ws.Cells.Item(4, 6).Value = "2021-10-05 13:39:34.561586"
$ws.Cells.Item(5, 6).Value = "2021-10-05 13:39:34.561589"
$ws.Cells.Item(6, 6).Value = "2021-10-05 13:39:34.561592"
$ws.Cells.Item(7, 6).Value = "2021-10-05 13:39:34.561596"
$ws.Cells.Item(8, 6).Value = "2021-10-05 13:39:34.561598"
$ws.Cells.Item(9, 6).Value = "2021-10-05 13:39:34.561601"
$ws.Cells.Item(10, 6).Value = "2021-10-05 13:39:34.561605"
$ws.Cells.Item(11, 6).Value = "2021-10-05 13:39:34.561608"
$ws.Cells.Item(12, 6).Value = "2021-10-05 13:39:34.561611"
$ws.Cells.Item(13, 6).Value = "2021-10-05 13:39:34.561614"
$ws.Cells.Item(14, 6).Value = "2021-10-05 13:39:34.561617"
$ws.Cells.Item(15, 6).Value = "2021-10-05 13:39:34.561619"
$ws.Cells.Item(16, 6).Value = "2021-10-05 13:39:34.561622"
$ws.Cells.Item(17, 6).Value = "2021-10-05 13:39:34.561625"
$ws.Cells.Item(18, 6).Value = "2021-10-05 13:39:34.561629"
$ws.Cells.Item(19, 6).Value = "2021-10-05 13:39:34.561632"
$ws.Cells.Item(20, 6).Value = "2021-10-05 13:39:34.561635"
$ws.Cells.Item(21, 6).Value = "2021-10-05 13:39:34.561638"
$ws.Cells.Item(22, 6).Value = "2021-10-05 13:39:34.561641"
$ws.Cells.Item(23, 6).Value = "2021-10-05 13:39:34.561643"
$ws.Cells.Item(24, 6).Value = "2021-10-05 13:39:34.561646"
$ws.Cells.Item(25, 6).Value = "2021-10-05 13:39:34.561649"
$ws.Cells.Item(26, 6).Value = "2021-10-05 13:39:34.561652"
$ws.Cells.Item(27, 6).Value = "2021-10-05 13:39:34.561655"
$ws.Cells.Item(28, 6).Value = "2021-10-05 13:39:34.561658"
$ws.Cells.Item(29, 6).Value = "2021-10-05 13:39:34.561662"
$ws.Cells.Item(30, 6).Value = "2021-10-05 13:39:34.561665"
$ws.Cells.Item(31, 6).Value = "2021-10-05 13:39:34.561668"
$ws.Cells.Item(32, 6).Value = "2021-10-05 13:39:34.561671"
$ws.Cells.Item(33, 6).Value = "2021-10-05 13:39:34.561674"
$ws.Cells.Item(34, 6).Value = "2021-10-05 13:39:34.561677"
$ws.Cells.Item(35, 6).Value = "2021-10-05 13:39:34.561680"
$ws.Cells.Item(36, 6).Value = "2021-10-05 13:39:34.561683"
$ws.Cells.Item(37, 6).Value = "2021-10-05 13:39:34.561686"
$ws.Cells.Item(38, 6).Value = "2021-10-05 13:39:34.561688"
$ws.Cells.Item(39, 6).Value = "2021-10-05 13:39:34.561691"
$ws.Cells.Item(40, 6).Value = "2021-10-05 13:39:34.561694"
$ws.Cells.Item(41, 6).Value = "2021-10-05 13:39:34.561697"
$ws.Cells.Item(42, 6).Value = "2021-10-05 13:39:34.561701"
$ws.Cells.Item(43, 6).Value = "2021-10-05 13:39:34.561704"
$ws.Cells.Item(44, 6).Value = "2021-10-05 13:39:34.561707"
$ws.Cells.Item(45, 6).Value = "2021-10-05 13:39:34.561710"
$ws.Cells.Item(46, 6).Value = "2021-10-05 13:39:34.561713"
$ws.Cells.Item(47, 6).Value = "2021-10-05 13:39:34.561716"
$ws.Cells.Item(48, 6).Value = "2021-10-05 13:39:34.561718"
$ws.Cells.Item(49, 6).Value = "2021-10-05 13:39:34.561721"
$ws.Cells.Item(50, 6).Value = "2021-10-05 13:39:34.561724"
$ws.Cells.Item(51, 6).Value = "2021-10-05 13:39:34.561727"
$ws.Cells.Item(52, 6).Value = "2021-10-05 13:39:34.561730"
$ws.Cells.Item(53, 6).Value = "2021-10-05 13:39:34.561733"
$ws.Cells.Item(54, 6).Value = "2021-10-05 13:39:34.561736"
$ws.Cells.Item(55, 6).Value = "2021-10-05 13:39:34.561739"
$ws.Cells.Item(56, 6).Value = "2021-10-05 13:39:34.561742"
$ws.Cells.Item(57, 6).Value = "2021-10-05 13:39:34.561745"
$ws.Cells.Item(58, 6).Value = "2021-10-05 13:39:34.561748"
$ws.Cells.Item(59, 6).Value = "2021-10-05 13:39:34.561751"
$ws.Cells.Item(60, 6).Value = "2021-10-05 13:39:34.561754"
$ws.Cells.Item(61, 6).Value = "2021-10-05 13:39:34.561757"
$ws.Cells.Item(62, 6).Value = "2021-10-05 13:39:34.561760"
$ws.Cells.Item(63, 6).Value = "2021-10-05 13:39:34.561763"
$ws.Cells.Item(64, 6).Value = "2021-10-05 13:39:34.561766"
$ws.Cells.Item(65, 6).Value = "2021-10-05 13:39:34.561769"
$ws.Cells.Item(66, 6).Value = "2021-10-05 13:39:34.561773"
$ws.Cells.Item(67, 6).Value = "2021-10-05 13:39:34.561776"
$ws.Cells.Item(68, 6).Value = "2021-10-05 13:39:34.561779"
$ws.Cells.Item(69, 6).Value = "2021-10-05 13:39:34.561782"
$ws.Cells.Item(70, 6).Value = "2021-10-05 13:39:34.561785"
$ws.Cells.Item(71, 6).Value = "2021-10-05 13:39:34.561788"
$ws.Cells.Item(72, 6).Value = "2021-10-05 13:39:34.561791"
$ws.Cells.Item(73, 6).Value = "2021-10-05 13:39:34.561794"
$ws.Cells.Item(74, 6).Value = "2021-10-05 13:39:34.561797"
$ws.Cells.Item(75, 6).Value = "2021-10-05 13:39:34.561800"
$ws.Cells.Item(76, 6).Value = "2021-10-05 13:39:34.561802"
$ws.Cells.Item(77, 6).Value = "2021-10-05 13:39:34.561805"
$ws.Cells.Item(78, 6).Value = "2021-10-05 13:39:34.561810"
$ws.Cells.Item(79, 6).Value = "2021-10-05 13:39:34.561814"
$ws.Cells.Item(80, 6).Value = "2021-10-05 13:39:34.561817"
$ws.Cells.Item(81, 6).Value = "2021-10-05 13:39:34.561820"

$wb.Save()
